$d = $word.ActiveDocument

# The document has three inline pictures, all living in headers/footers:
#   - header1.xml ("first" header)   -> BTec logo,    docPr/cNvPr name "image1.jpg" -> "image2.jpg"
#   - footer1.xml ("default" footer) -> Pearson logo,  docPr/cNvPr name "image2.png" -> "image1.png"
#   - footer2.xml ("first" footer)   -> Pearson logo,  docPr/cNvPr name "image2.png" -> "image1.png"
#
# InlineShape has no writable Name property (matches real Word), so each
# picture is round-tripped through Shape (ConvertToShape -> rename ->
# ConvertToInlineShape) which updates the shape's docPr/name while leaving
# it as a plain <wp:inline> picture, same as before.

function Rename-InlinePicture($range, [string]$newName) {
    $shapeCount = $range.InlineShapes.Count
    for ($i = 1; $i -le $shapeCount; $i++) {
        $inlineShape = $range.InlineShapes.Item($i)
        $shape = $inlineShape.ConvertToShape()
        $shape.Name = $newName
        [void]$shape.ConvertToInlineShape()
    }
}

$sec = $d.Sections(1)

# Header (type "first"): BTec_Logo-Orange, image1.jpg -> image2.jpg
for ($h = 1; $h -le $sec.Headers.Count; $h++) {
    $hf = $sec.Headers($h)
    if ($hf.Exists -and $hf.Range.InlineShapes.Count -gt 0) {
        Rename-InlinePicture $hf.Range "image2.jpg"
    }
}

# Footers (both "default" and "first"): PearsonLogo, image2.png -> image1.png
for ($f = 1; $f -le $sec.Footers.Count; $f++) {
    $ft = $sec.Footers($f)
    if ($ft.Exists -and $ft.Range.InlineShapes.Count -gt 0) {
        Rename-InlinePicture $ft.Range "image1.png"
    }
}
